$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns for affected rows so Excel does not
# reinterpret values like "3.51" or "1.00" as numbers, or "65.555.66" etc.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.555.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.273.52'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.97%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.27'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.54%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.261.41'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.83%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.587'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.49'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -8.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000264'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '634.03'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.55'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.814.88'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.667.63'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.85'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.116'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.284.98'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -5.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.36'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -8.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.905'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.84'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '107.01'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.99'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -6.98%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -7.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.68'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.50'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.69'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.30'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.96'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.33'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.06'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '555.43'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +11.15%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.14'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.81%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.654.75'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.85'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +43.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.49'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.72'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0711'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -9.39%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.342'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.03'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0416'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.24'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.129'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.66%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.24'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.50%  '
